$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format D-column cells whose new price value would otherwise be
# auto-converted from text to a number (e.g. "1.00" -> 1) by Excel,
# so they keep their exact textual representation, as in the source data.
$textCells = @("D4","D5","D6","D7","D8","D10","D11","D14","D15","D16","D18","D20","D21","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin values (price, volume %, and for three rows
# whose rank order changed, the coin name + link as well).
$ws.Range("D2").Value = '27.942.62'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '1.638.84'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '212.65'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").Value = '0.523'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '23.42'
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("D10").Value = '0.0613'
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = '0.0883'
$ws.Range("E11").Value = '  +2.41%  '
$ws.Range("D12").Value = '1.870.65'
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").Value = '1.638.31'
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").Value = '4.08'
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("D15").Value = '0.572'
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").Value = '65.47'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '27.925.81'
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").Value = '232.48'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").Value = '0.0₃0721'
$ws.Range("D20").Value = '7.59'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '10.46'
$ws.Range("E22").Value = '  -2.62%  '
$ws.Range("D23").Value = '4.37'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '2.09'
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").Value = '152.67'
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").Value = '6.89'
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").Value = '15.69'
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '1.19'
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("D31").Value = '0.0484'
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = '3.36'
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("D34").Value = '1.405.40'
$ws.Range("E34").Value = '  -4.09%  '
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  +2.28%  '
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("D37").Value = '0.0170'
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("D38").Value = '0.880'
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '0.926'
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").Value = '0.559'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").Value = '1.03'
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").Value = '67.36'
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '1.85'
$ws.Range("E44").Value = '  +6.18%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '5.51'
$ws.Range("E45").Value = '  +2.57%  '
$ws.Range("B46").Value = 'MXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D46").Value = '2.21'
$ws.Range("E46").Value = '  -2.83%  '
$ws.Range("D47").Value = '1.779.71'
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").Value = '87.72'
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("E49").Value = '  +0.55%  '
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = '7.64'
$ws.Range("E51").Value = '  -0.24%  '
